# "La till ett !" - append an extra "!" to the end of the document's last
# paragraph, which currently reads "Och detta är INTE en låt!" so that it
# becomes "Och detta är INTE en låt!!".
#
# The tricky part is that Word's hidden "_GoBack" bookmark (which marks the
# location of the most recent edit) originally sits in the middle of the
# last paragraph (between "Och d" and "etta är INTE en låt!") and needs to
# move to the very end of the paragraph, right after the text we add,
# mirroring what real Word does when you click at the end of the document
# and type a character.
#
# Note: this COM-interop runtime has a quirk where Bookmarks.Add() with a
# *collapsed* Range positioned exactly at the last character slot of a
# paragraph (i.e. immediately before the paragraph mark) silently produces
# a corrupted bookmark range. We work around it by temporarily inserting
# one extra placeholder character after the real insertion point, adding
# the bookmark while it is safely *not* at the paragraph-end boundary, and
# then deleting the placeholder - the bookmark (being to the left of the
# deleted text) stays put at the correct, now-final position.

$d = $word.ActiveDocument

# Locate the last paragraph and a collapsed Range sitting right at the end
# of its visible text, i.e. immediately before the paragraph mark.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertRange = $lastPara.Range.Duplicate()
$insertRange.Collapse(0)        # 0 = wdCollapseEnd
[void]$insertRange.MoveEnd(1, -1)     # step back over the paragraph mark
$insertPos = $insertRange.Start

# Type the new "!" - plus one temporary placeholder "!" that lets us place
# the bookmark safely before removing it again.
$insertRange.InsertAfter("!!")

# The real new character now sits at $insertPos, the placeholder right
# after it at $insertPos + 1. The bookmark belongs right between them,
# i.e. right after the genuine new "!".
$bmPos = $insertPos + 1
$bmRange = $d.Range($bmPos, $bmPos)

# Move "_GoBack" there: drop the old one (if present) and add it fresh.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary placeholder character; the bookmark (anchored just
# to its left) is unaffected and ends up correctly collapsed at the new
# end of the paragraph.
$d.Range($bmPos, $bmPos + 1).Delete()
